$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the status text that is shared across Overview / zh-cn / de-de
#    ("In Translation" -> "Handed back: in sync with en-US")
# ------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("In Translation", "Handed back: in sync with en-US", 1, 1, $false, $false, $false, $false) | Out-Null
}

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8432e3fd9b6663b5cfd038c575d8f8641db3b84b/e2e/824f3a10-d14e-4400-97fc-25b74f189dfc.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8432e3fd9b6663b5cfd038c575d8f8641db3b84b/e2e/96b64975-f652-40e6-81ac-ecc76a189561.md"
$mdName1 = "824f3a10-d14e-4400-97fc-25b74f189dfc.md"
$mdName2 = "96b64975-f652-40e6-81ac-ecc76a189561.md"

# ------------------------------------------------------------------
# 2. zh-cn sheet: fill in the Latest Target File / Latest Handback File /
#    Latest Handback DateTime columns for the two rows and add the new
#    hyperlinks on the Latest Target File column (I).
# ------------------------------------------------------------------
$zh.Range("I2").Value = $mdName1
$zh.Range("J2").Value = "824f3a10-d14e-4400-97fc-25b74f189dfc.dfb1cb8d4e03c17f255cfcba5c992f802dc123ef.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-30 09:32:22"

$zh.Range("I3").Value = $mdName2
$zh.Range("J3").Value = "96b64975-f652-40e6-81ac-ecc76a189561.f829a9294e04b4626ea0a925866b60d27730a2bf.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-30 09:32:22"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $mdUrl1, "", "", $mdName1)
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl1, "", "", $mdName1)
$zh.Hyperlinks.Add($zh.Range("A3"), $mdUrl2, "", "", $mdName2)
$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl2, "", "", $mdName2)

# ------------------------------------------------------------------
# 3. de-de sheet: same treatment, with its own handback timestamp.
# ------------------------------------------------------------------
$de.Range("I2").Value = $mdName1
$de.Range("J2").Value = "824f3a10-d14e-4400-97fc-25b74f189dfc.dfb1cb8d4e03c17f255cfcba5c992f802dc123ef.de-de.xlf"
$de.Range("K2").Value = "2016-08-30 09:32:42"

$de.Range("I3").Value = $mdName2
$de.Range("J3").Value = "96b64975-f652-40e6-81ac-ecc76a189561.f829a9294e04b4626ea0a925866b60d27730a2bf.de-de.xlf"
$de.Range("K3").Value = "2016-08-30 09:32:42"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $mdUrl1, "", "", $mdName1)
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl1, "", "", $mdName1)
$de.Hyperlinks.Add($de.Range("A3"), $mdUrl2, "", "", $mdName2)
$de.Hyperlinks.Add($de.Range("I3"), $mdUrl2, "", "", $mdName2)

# ------------------------------------------------------------------
# 4. Column widths grow because of the longer text that now lives in
#    them (status text, guid filenames, xlf filenames) - mimic the
#    resulting auto-fit sizes.
# ------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
